$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit permutes the full content of data rows 2-9 and 12-24 among
# themselves (rows 10 and 11 are left untouched). Each row's entire set of
# cells (values, typed-but-empty placeholder cells, booleans, etc.) moves as
# a unit to a different row position, mirroring the upstream diff.
#
# The permutation decomposes into independent cycles; each cycle is rotated
# using a temporary holding row (far outside the used range) so that no data
# is lost while values are shuffled around.

function Rotate-Cycle {
    param(
        [int[]]$Cycle,
        $Worksheet,
        [int]$TempRow
    )

    $first = $Cycle[0]
    $firstRange = "A" + $first + ":AY" + $first
    $tempRange = "A" + $TempRow + ":AY" + $TempRow
    $Worksheet.Range($firstRange).Copy($Worksheet.Range($tempRange))

    # Excel's Range.Copy only writes cells that are present in the source
    # range; it does not blank out destination cells that have no
    # corresponding source cell. Since some rows carry "extra" placeholder
    # cells (e.g. typed-but-empty cells, or substrate columns) that must not
    # survive the move, clear every destination row before pasting into it.
    for ($i = 0; $i -lt ($Cycle.Count - 1); $i++) {
        $destRow = $Cycle[$i]
        $srcRow = $Cycle[$i + 1]
        $destRange = "A" + $destRow + ":AY" + $destRow
        $srcRange = "A" + $srcRow + ":AY" + $srcRow
        $Worksheet.Range($destRange).ClearContents()
        $Worksheet.Range($srcRange).Copy($Worksheet.Range($destRange))
    }

    $last = $Cycle[$Cycle.Count - 1]
    $lastRange = "A" + $last + ":AY" + $last
    $Worksheet.Range($lastRange).ClearContents()
    $Worksheet.Range($tempRange).Copy($Worksheet.Range($lastRange))
    $Worksheet.Range($tempRange).ClearContents()
}

$cycles = @(
    @(2, 6, 3, 5, 4),
    @(7, 23, 21, 16, 9, 17, 12, 14, 18),
    @(8, 20, 22),
    @(13, 19, 15, 24)
)

foreach ($cycle in $cycles) {
    Rotate-Cycle $cycle $ws 1000
}
